$wb = $excel.ActiveWorkbook
$contacts = $wb.Worksheets.Item("contacts")

# Remove the stray D34 cell (value was a single space) from the contacts sheet
$contacts.Range("D34").ClearContents()

# Add a new sheet with a copy of the fname/lname table, placed after "contacts"
$newSheet = $wb.Worksheets.Add($null, $contacts)
$newSheet.Name = "Sheet1"
for ($r = 1; $r -le 5; $r++) {
  $newSheet.Cells.Item($r, 1).Value = $contacts.Cells.Item($r, 1).Value2
  $newSheet.Cells.Item($r, 2).Value = $contacts.Cells.Item($r, 2).Value2
}
# Match the yellow header-row highlight used on "contacts"
$newSheet.Range("A1:B1").Interior.Color = 65535

# Set the new sheet's selection (whole table selected, active cell at B5)
$newSheet.Range("A1:B5").Select()
$newSheet.Range("B5").Activate()

# Re-select contacts as the active/visible tab, with A1:B5 selected (active cell A1)
$contacts.Select()
$contacts.Range("A1:B5").Select()

$wb.Save()
